$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$v = $ws.Range("B2").Value()
"value is: $v"
